# The commit adds a "result" (PASS/FAIL) column next to the existing
# username/password test data on the first worksheet, reflecting the
# outcome of each login test case once the data-provider was split out
# into its own class.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "FAIL"
$ws.Range("C3").Value = "PASS"
$ws.Range("C4").Value = "FAIL"
$ws.Range("C5").Value = "FAIL"
